$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.225.97'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.28%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.772.00'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.52%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.70'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.35%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5249'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +10.33%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3667'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +6.66%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '42.66'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.29%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07350'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.02%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.086'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.24%  '

$ws.Range("E12").Value = '  -0.02%  '

$ws.Range("E13").Value = '  +2.93%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.052'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.31%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.769.10'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.38%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.922'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.13%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '88.62'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.13%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001043'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.21%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06433'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.21%  '

$ws.Range("E20").Value = '  -0.01%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.69'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.25%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.802'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.48%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '27.269.84'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.32%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.24'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.06%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.126'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.67%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '155.07'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.10%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.08'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.10%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.972.98'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.51%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.322'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +11.65%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '120.97'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.83%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.055'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.98%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09775'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.67%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.551'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.86%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.621'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.83%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02231'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.80%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05958'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.11%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '11.19'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.31%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2017'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.36%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.819'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.46%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6117'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.33%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.432'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.37%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.060'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +7.74%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.137'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.11%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.14'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.88%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5751'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.33%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.620'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.52%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '121.17'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.47%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.872'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.80%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.113'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.59%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06698'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.96%  '

$ws.Range("E51").Value = '  +0.02%  '
